# Apply the Mandragora_Profits scheduled-runner update: refreshed market-board
# price snapshots (currentAveragePrice* / LevePrice* / LeveProfit* columns)
# across the per-job leve sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 17
$ws.Range("H17").Value = 956.51514
$ws.Range("J17").Value = 956.51514
$ws.Range("L17").Value = 2869.54542
$ws.Range("N17").Value = -3205.54542
# row 19
$ws.Range("H19").Value = 50000944
$ws.Range("I19").Value = 599.3333
$ws.Range("J19").Value = 71429660
$ws.Range("K19").Value = 599.3333
$ws.Range("L19").Value = 71429660
$ws.Range("M19").Value = -424.3333
$ws.Range("N19").Value = -71430010
# row 100
$ws.Range("H100").Value = 2281.818
$ws.Range("I100").Value = 2180
$ws.Range("J100").Value = 2366.6667
$ws.Range("K100").Value = 2180
$ws.Range("L100").Value = 2366.6667
$ws.Range("M100").Value = -1639
$ws.Range("N100").Value = -3448.6667
# row 127
$ws.Range("H127").Value = 1146.2632
$ws.Range("J127").Value = 1357
$ws.Range("L127").Value = 4071
$ws.Range("N127").Value = -13991
# row 129
$ws.Range("H129").Value = 1529.3334
$ws.Range("I129").Value = 767.55554
$ws.Range("J129").Value = 1719.7778
$ws.Range("K129").Value = 2302.66662
$ws.Range("L129").Value = 5159.3334
$ws.Range("M129").Value = 2697.33338
$ws.Range("N129").Value = -15159.3334
# row 137
$ws.Range("H137").Value = 19563.691
$ws.Range("I137").Value = 1917.6957
$ws.Range("J137").Value = 32246.75
$ws.Range("K137").Value = 5753.0871
$ws.Range("L137").Value = 96740.25
$ws.Range("M137").Value = -3203.0871
$ws.Range("N137").Value = -101840.25

$ws = $wb.Worksheets.Item("ARM")
# row 57
$ws.Range("H57").Value = 26000
$ws.Range("I57").Value = 26000
$ws.Range("K57").Value = 26000
$ws.Range("M57").Value = -25516
# row 61
$ws.Range("H61").Value = 1874.7556
$ws.Range("I61").Value = 1881.75
$ws.Range("J61").Value = 1863.2354
$ws.Range("K61").Value = 1881.75
$ws.Range("L61").Value = 1863.2354
$ws.Range("M61").Value = -1669.75
$ws.Range("N61").Value = -2287.2354
# row 62
$ws.Range("H62").Value = 20000
$ws.Range("J62").Value = 20000
$ws.Range("L62").Value = 20000
$ws.Range("N62").Value = -21248
# row 65
$ws.Range("H65").Value = 20000
$ws.Range("J65").Value = 20000
$ws.Range("L65").Value = 60000
$ws.Range("N65").Value = -66240
# row 126
$ws.Range("H126").Value = 5333.3335
$ws.Range("I126").Value = 5333.3335
$ws.Range("K126").Value = 16000.0005
$ws.Range("M126").Value = -13530.0005
# row 132
$ws.Range("H132").Value = 5809.946
$ws.Range("I132").Value = 3231.2856
$ws.Range("J132").Value = 9194.4375
$ws.Range("K132").Value = 9693.856800000001
$ws.Range("L132").Value = 27583.3125
$ws.Range("M132").Value = -7163.856800000001
$ws.Range("N132").Value = -32643.3125
# row 133
$ws.Range("H133").Value = 40000
$ws.Range("J133").Value = 40000
$ws.Range("L133").Value = 40000
$ws.Range("N133").Value = -45060
# row 136
$ws.Range("H136").Value = 1874.7556
$ws.Range("I136").Value = 1881.75
$ws.Range("J136").Value = 1863.2354
$ws.Range("K136").Value = 5645.25
$ws.Range("L136").Value = 5589.706200000001
$ws.Range("M136").Value = -3095.25
$ws.Range("N136").Value = -10689.7062

$ws = $wb.Worksheets.Item("BSM")
# row 134
$ws.Range("H134").Value = 2863.8823
$ws.Range("I134").Value = 1653.2273
$ws.Range("K134").Value = 4959.6819
$ws.Range("M134").Value = -2424.6819

$ws = $wb.Worksheets.Item("CRP")
# row 69
$ws.Range("H69").Value = 8300
$ws.Range("I69").Value = 8300
$ws.Range("K69").Value = 8300
$ws.Range("M69").Value = -7551
# row 72
$ws.Range("H72").Value = 8300
$ws.Range("I72").Value = 8300
$ws.Range("K72").Value = 24900
$ws.Range("M72").Value = -21156

$ws = $wb.Worksheets.Item("CUL")
# row 5
$ws.Range("H5").Value = 2565
$ws.Range("I5").Value = 2570
$ws.Range("J5").Value = 2561.6667
$ws.Range("K5").Value = 7710
$ws.Range("L5").Value = 7685.000100000001
$ws.Range("M5").Value = -7598
$ws.Range("N5").Value = -7909.000100000001
# row 36
$ws.Range("H36").Value = 2963.9285
$ws.Range("J36").Value = 4362.875
$ws.Range("L36").Value = 13088.625
$ws.Range("N36").Value = -13426.625
# row 46
$ws.Range("H46").Value = 635.6667
$ws.Range("I46").Value = 303
$ws.Range("K46").Value = 909
$ws.Range("M46").Value = -818
# row 60
$ws.Range("H60").Value = 1130
$ws.Range("I60").Value = 273.75
$ws.Range("J60").Value = 2500
$ws.Range("K60").Value = 821.25
$ws.Range("L60").Value = 7500
$ws.Range("M60").Value = -570.25
$ws.Range("N60").Value = -8002
# row 135
$ws.Range("H135").Value = 2565
$ws.Range("I135").Value = 2570
$ws.Range("J135").Value = 2561.6667
$ws.Range("K135").Value = 23130
$ws.Range("L135").Value = 23055.0003
$ws.Range("M135").Value = -20595
$ws.Range("N135").Value = -28125.0003

$ws = $wb.Worksheets.Item("GSM")
# row 80
$ws.Range("H80").Value = 2434
$ws.Range("I80").Value = 2461
$ws.Range("J80").Value = 2414.7144
$ws.Range("K80").Value = 2461
$ws.Range("L80").Value = 2414.7144
$ws.Range("M80").Value = -1463
$ws.Range("N80").Value = -4410.7144
# row 83
$ws.Range("H83").Value = 2434
$ws.Range("I83").Value = 2461
$ws.Range("J83").Value = 2414.7144
$ws.Range("K83").Value = 12305
$ws.Range("L83").Value = 12073.572
$ws.Range("M83").Value = -7313
$ws.Range("N83").Value = -22057.572
# row 93
$ws.Range("H93").Value = 19750
$ws.Range("J93").Value = 19750
$ws.Range("L93").Value = 19750
$ws.Range("N93").Value = -23494

$ws = $wb.Worksheets.Item("LTW")
# row 62
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
# row 65
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
# row 93
$ws.Range("H93").Value = 2489.8948
$ws.Range("I93").Value = 2314.2856
$ws.Range("J93").Value = 2981.6
$ws.Range("K93").Value = 2314.2856
$ws.Range("L93").Value = 2981.6
$ws.Range("M93").Value = -1066.2856
$ws.Range("N93").Value = -5477.6
# row 132
$ws.Range("H132").Value = 26110.795
$ws.Range("I132").Value = 32737.295
$ws.Range("J132").Value = 3580.7
$ws.Range("K132").Value = 98211.88499999999
$ws.Range("L132").Value = 10742.1
$ws.Range("M132").Value = -95681.88499999999
$ws.Range("N132").Value = -15802.1
# row 136
$ws.Range("H136").Value = 2369.077
$ws.Range("I136").Value = 2379.8
$ws.Range("J136").Value = 2333.3333
$ws.Range("K136").Value = 7139.400000000001
$ws.Range("L136").Value = 6999.999899999999
$ws.Range("M136").Value = -4589.400000000001
$ws.Range("N136").Value = -12099.9999

$ws = $wb.Worksheets.Item("WVR")
# row 82
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
# row 85
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
# row 136
$ws.Range("H136").Value = 1496.3143
$ws.Range("I136").Value = 1153.7637
$ws.Range("J136").Value = 2752.3333
$ws.Range("K136").Value = 3461.2911
$ws.Range("L136").Value = 8256.999899999999
$ws.Range("M136").Value = -911.2910999999999
$ws.Range("N136").Value = -13356.9999
